$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gas_criteria")

$ws.Range("B2").Value = 971
$ws.Range("C2").Value = 49

$ws.Range("B3").Value = 1300
$ws.Range("C3").Value = 65.59999999999999

$ws.Range("B4").Value = 330
$ws.Range("C4").Value = 16.6

$ws.Range("B5").Value = 1523
$ws.Range("C5").Value = 76.8
